# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Handback timestamps for zh-cn / de-de are refreshed
#  - The stale "handback file is not latest" error is cleared
#  - Column widths are widened/narrowed to fit the new text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status summary columns ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-16 00:45:17"
$zhcn.Range("P2").Value = ""

# --- de-de sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-16 00:45:24"
$dede.Range("P2").Value = ""

# --- Column width adjustments to fit the new content ---
$overview.Range("E1").ColumnWidth = 29.166666666666668
$overview.Range("F1").ColumnWidth = 29.166666666666668

$zhcn.Range("C1").ColumnWidth = 29.166666666666668
$zhcn.Range("P1").ColumnWidth = 12.833333333333334

$dede.Range("C1").ColumnWidth = 29.166666666666668
$dede.Range("P1").ColumnWidth = 12.833333333333334
